$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value2 = $ws.Range("A6").Value2
$ws.Range("B7").Value2 = $ws.Range("B6").Value2
